$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (row 1)
$ws.Range("A1").Value = "p1"
$ws.Range("B1").Value = "p2"
$ws.Range("C1").Value = "p3"
$ws.Range("D1").Value = "p4"
$ws.Range("E1").Value = "p5"
$ws.Range("F1").Value = "p6"

# Add new row 2 values
$ws.Range("A2").Value = "no"
$ws.Range("B2").Value = "algunas veces"

# Materialize the remaining (empty) row-2 cells with the default style,
# matching the widened used-range (A1:F2) from the source sheet.
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Style = "Normal"

# Resize the columns to the new widths
$ws.Columns.Item(1).ColumnWidth = 4.571428571428571
$ws.Columns.Item(2).ColumnWidth = 10.142857142857142
$ws.Columns.Item(3).ColumnWidth = 4.571428571428571
$ws.Columns.Item(4).ColumnWidth = 4.571428571428571
$ws.Columns.Item(5).ColumnWidth = 4.571428571428571
$ws.Columns.Item(6).ColumnWidth = 4.571428571428571
